$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add summary rows 27 (MIN) and 28 (MAX) for columns B, F, G ---
$ws.Range("B27").Formula = "=MIN(B3,B5:B8,B10:B21,B23,B25)"
$ws.Range("F27").Formula = "=MIN(F2:F25)"
$ws.Range("G27").Formula = "=MIN(G2:G25)"

$ws.Range("B28").Formula = "=MAX(B2:B25)"
$ws.Range("F28").Formula = "=MAX(F2:F25)"
$ws.Range("G28").Formula = "=MAX(G2:G25)"

# --- Turn the URL text already in H2 into a real hyperlink ---
$ws.Hyperlinks.Add($ws.Range("H2"), "https://www.cscu.cornell.edu/news/Handouts/SEM_fit.pdf") | Out-Null

# --- Move the active selection to the new last cell, like the author did ---
$ws.Range("H27").Select() | Out-Null
